# Generate Report for Handoff
# Stamp fresh "handoff" timestamps for the b812e9d9-5fd8-44da-ba73-9a47a60d8461
# report row (row 7) across the Overview, zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-09-06 08:55:58"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-09-06 08:55:53"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-09-06 08:55:58"
